# Update the date label and regenerate the practice-problem answers.
# Each old answer string is unique in the document, so a single
# Find/Replace (ReplaceAll=2) per pair is unambiguous. The pairs are
# applied in document order, which also guarantees that the lone
# old/new text collision ("19÷4=4, 3" is both an original value and a
# freshly produced value later on) is resolved safely: the original
# occurrence is renamed away before the new occurrence is introduced.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-09 Sunday" "2025-03-10 Monday"

Replace-Text "78÷6=13, 0" "39÷3=13, 0"
Replace-Text "25÷6=4, 1" "65÷5=13, 0"
Replace-Text "67÷3=22, 1" "23÷9=2, 5"
Replace-Text "47÷2=23, 1" "94÷6=15, 4"
Replace-Text "29÷6=4, 5" "91÷9=10, 1"

Replace-Text "15÷5=3, 0" "37÷2=18, 1"
Replace-Text "87÷4=21, 3" "62÷8=7, 6"
Replace-Text "88÷7=12, 4" "79÷9=8, 7"
Replace-Text "14÷6=2, 2" "14÷7=2, 0"
Replace-Text "19÷4=4, 3" "74÷4=18, 2"

Replace-Text "30÷6=5, 0" "45÷7=6, 3"
Replace-Text "95÷3=31, 2" "56÷6=9, 2"
Replace-Text "19÷9=2, 1" "50÷8=6, 2"
Replace-Text "40÷7=5, 5" "50÷4=12, 2"
Replace-Text "52÷5=10, 2" "59÷8=7, 3"

Replace-Text "14÷5=2, 4" "55÷4=13, 3"
Replace-Text "83÷6=13, 5" "58÷2=29, 0"
Replace-Text "48÷8=6, 0" "19÷2=9, 1"
Replace-Text "70÷9=7, 7" "98÷6=16, 2"
Replace-Text "37÷3=12, 1" "72÷8=9, 0"

Replace-Text "16÷2=8, 0" "95÷8=11, 7"
Replace-Text "81÷6=13, 3" "43÷8=5, 3"
Replace-Text "80÷4=20, 0" "19÷4=4, 3"
Replace-Text "25÷3=8, 1" "83÷4=20, 3"
Replace-Text "92÷6=15, 2" "75÷5=15, 0"

Write-Output "replacements complete"
